$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H39").Value = 447.8
$ws.Range("I39").Value = 113.333336
$ws.Range("J39").Value = 949.5
$ws.Range("K39").Value = 340.000008
$ws.Range("L39").Value = 2848.5
$ws.Range("M39").Value = -44.00000799999998
$ws.Range("N39").Value = -3440.5

$ws.Range("H64").Value = 14096487
$ws.Range("I64").Value = 3658308.2
$ws.Range("J64").Value = 35718428
$ws.Range("K64").Value = 3658308.2
$ws.Range("L64").Value = 35718428
$ws.Range("M64").Value = -3658060.2
$ws.Range("N64").Value = -35718924

$ws.Range("H67").Value = 14096487
$ws.Range("I67").Value = 3658308.2
$ws.Range("J67").Value = 35718428
$ws.Range("K67").Value = 3658308.2
$ws.Range("L67").Value = 35718428
$ws.Range("M67").Value = -3657450.2
$ws.Range("N67").Value = -35720144

$ws.Range("H82").Value = 5147.8887
$ws.Range("I82").Value = 1266.2
$ws.Range("K82").Value = 3798.6
$ws.Range("M82").Value = -3392.6

$ws.Range("H85").Value = 5147.8887
$ws.Range("I85").Value = 1266.2
$ws.Range("K85").Value = 3798.6
$ws.Range("M85").Value = -2394.6

$ws.Range("H129").Value = 1457.6923
$ws.Range("I129").Value = 844.375
$ws.Range("J129").Value = 2439
$ws.Range("K129").Value = 2533.125
$ws.Range("L129").Value = 7317
$ws.Range("M129").Value = 2466.875
$ws.Range("N129").Value = -17317

$ws.Range("H136").Value = 153166.67
$ws.Range("J136").Value = 153166.67
$ws.Range("L136").Value = 153166.67
$ws.Range("N136").Value = -163366.67

$ws.Range("H138").Value = 3130.65
$ws.Range("J138").Value = 3528.125
$ws.Range("L138").Value = 10584.375
$ws.Range("N138").Value = -20864.375

$ws.Range("H140").Value = 66751
$ws.Range("J140").Value = 66185.57000000001
$ws.Range("L140").Value = 66185.57000000001
$ws.Range("N140").Value = -76545.57000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1343892.4
$ws.Range("J2").Value = 2439.8333
$ws.Range("L2").Value = 2439.8333
$ws.Range("N2").Value = -2665.8333

$ws.Range("H32").Value = 4141.1313
$ws.Range("I32").Value = 2060.48
$ws.Range("J32").Value = 13598.637
$ws.Range("K32").Value = 2060.48
$ws.Range("L32").Value = 13598.637
$ws.Range("M32").Value = -1773.48
$ws.Range("N32").Value = -14172.637

$ws.Range("H61").Value = 15301.454
$ws.Range("I61").Value = 15301.454
$ws.Range("K61").Value = 15301.454
$ws.Range("M61").Value = -15089.454

$ws.Range("H96").Value = 0
$ws.Range("J96").Value = 0
$ws.Range("N96").ClearContents()

$ws.Range("H116").Value = 1343892.4
$ws.Range("J116").Value = 2439.8333
$ws.Range("L116").Value = 2439.8333
$ws.Range("N116").Value = -7027.8333

$ws.Range("H132").Value = 5583.9697
$ws.Range("I132").Value = 2101.5
$ws.Range("J132").Value = 8861.588
$ws.Range("K132").Value = 6304.5
$ws.Range("L132").Value = 26584.764
$ws.Range("M132").Value = -3774.5
$ws.Range("N132").Value = -31644.764

$ws.Range("H136").Value = 15301.454
$ws.Range("I136").Value = 15301.454
$ws.Range("K136").Value = 45904.362
$ws.Range("M136").Value = -43354.362

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1343892.4
$ws.Range("J3").Value = 2439.8333
$ws.Range("L3").Value = 2439.8333
$ws.Range("N3").Value = -2667.8333

$ws.Range("H20").Value = 2883.3572
$ws.Range("I20").Value = 3113.0833
$ws.Range("J20").Value = 1505
$ws.Range("K20").Value = 3113.0833
$ws.Range("L20").Value = 1505
$ws.Range("M20").Value = -2866.0833
$ws.Range("N20").Value = -1999

$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 0
$ws.Range("K26").Value = 0

$ws.Range("H134").Value = 4094
$ws.Range("I134").Value = 1270
$ws.Range("J134").Value = 8330
$ws.Range("K134").Value = 3810
$ws.Range("L134").Value = 24990
$ws.Range("M134").Value = -1275
$ws.Range("N134").Value = -30060

$ws.Range("H140").Value = 99994
$ws.Range("J140").Value = 99994
$ws.Range("L140").Value = 99994
$ws.Range("N140").Value = -110354

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 265703
$ws.Range("I58").Value = 455165.38
$ws.Range("J58").Value = 5192.25
$ws.Range("K58").Value = 455165.38
$ws.Range("L58").Value = 5192.25
$ws.Range("M58").Value = -454962.38
$ws.Range("N58").Value = -5598.25

$ws.Range("H122").Value = 2853.4138
$ws.Range("I122").Value = 1725
$ws.Range("J122").Value = 5361
$ws.Range("K122").Value = 5175
$ws.Range("L122").Value = 16083
$ws.Range("M122").Value = -2725
$ws.Range("N122").Value = -20983

$ws.Range("H136").Value = 265703
$ws.Range("I136").Value = 455165.38
$ws.Range("J136").Value = 5192.25
$ws.Range("K136").Value = 1365496.14
$ws.Range("L136").Value = 15576.75
$ws.Range("M136").Value = -1362946.14
$ws.Range("N136").Value = -20676.75

$ws.Range("H141").Value = 107279.086
$ws.Range("J141").Value = 113095.45
$ws.Range("L141").Value = 113095.45
$ws.Range("N141").Value = -123455.45

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1482369.8
$ws.Range("I5").Value = 601.2857
$ws.Range("J5").Value = 2058613
$ws.Range("K5").Value = 1803.8571
$ws.Range("L5").Value = 6175839
$ws.Range("M5").Value = -1691.8571
$ws.Range("N5").Value = -6176063

$ws.Range("H14").Value = 305.72726
$ws.Range("I14").Value = 305.72726
$ws.Range("K14").Value = 917.18178
$ws.Range("M14").Value = -744.18178

$ws.Range("H46").Value = 2665.3333
$ws.Range("I46").Value = 999.75
$ws.Range("K46").Value = 2999.25
$ws.Range("M46").Value = -2908.25

$ws.Range("H121").Value = 250897.42
$ws.Range("J121").Value = 429330.84
$ws.Range("L121").Value = 1287992.52
$ws.Range("N121").Value = -1290612.52

$ws.Range("H122").Value = 5377055.5
$ws.Range("J122").Value = 899
$ws.Range("L122").Value = 8091
$ws.Range("N122").Value = -12991

$ws.Range("H126").Value = 16375
$ws.Range("J126").Value = 16833.334
$ws.Range("L126").Value = 50500.00199999999
$ws.Range("N126").Value = -60380.00199999999

$ws.Range("H129").Value = 4000
$ws.Range("I129").Value = 4000
$ws.Range("K129").Value = 12000
$ws.Range("M129").Value = -7000

$ws.Range("H131").Value = 8199272.5
$ws.Range("I131").Value = 13890942
$ws.Range("K131").Value = 41672826
$ws.Range("M131").Value = -41667786

$ws.Range("H135").Value = 1482369.8
$ws.Range("I135").Value = 601.2857
$ws.Range("J135").Value = 2058613
$ws.Range("K135").Value = 5411.571300000001
$ws.Range("L135").Value = 18527517
$ws.Range("M135").Value = -2876.571300000001
$ws.Range("N135").Value = -18532587

$ws.Range("H140").Value = 3649.9688
$ws.Range("I140").Value = 2779.4285
$ws.Range("K140").Value = 8338.2855
$ws.Range("M140").Value = -3158.2855

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H121").Value = 30317
$ws.Range("J121").Value = 30317
$ws.Range("L121").Value = 30317
$ws.Range("N121").Value = -33811
